# The upstream commit ("Moving from 2.0.2 to 2.0.3") only touches this
# template's serialized attribute ordering (namespace declarations on
# <w:document>, and attribute order on <w:pgSz>/<w:pgMar>, the
# <w:rFonts>/<w:lang> defaults, <w:latentStyles>/<w:lsdException>, and the
# four <w:style> definitions in styles.xml). Every changed line carries the
# exact same element, attributes and values as before - this is a
# round-trip re-serialization (the template was simply re-saved), not a
# content edit. There is nothing in the document's structure, text,
# formatting or properties that actually needs to change.
#
# Touch the document through the object model (equivalent to Word's
# "re-save" that produced the upstream diff) without altering any content.
$d = $word.ActiveDocument
$d.Save()
